# Applies the "Changed activity user screen" commit:
#  - Typography sheet: adds a new "Narrow" typography entry (row 7)
#  - Translation sheet: changes the "Activity" text id entry to "Cycling" (row 62)
#    and adds a new block of text ids for the activity user screen (rows 63-84)

$wb = $excel.ActiveWorkbook
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# Excel's COM layer eagerly parses plain .Value assignments that look like
# numbers (e.g. "0.00") into actual numeric cells. Several text ids in this
# sheet intentionally hold such number-shaped strings, so force those through
# as literal text via a text formula + paste-as-values round trip.
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# ---- Typography sheet: new row 7 ("Narrow" font) ----
$wsTypo.Range("B7").Value = "Narrow"
$wsTypo.Range("C7").Value = "ARIALN.TTF"
$wsTypo.Range("D7").Value = 28
$wsTypo.Range("E7").Value = 4
$wsTypo.Range("F7").Value = "?"
$wsTypo.Range("I7").Value = "0x00-0xFF"
# Columns B:E carry an explicit column style (Arial 10); writing .Value
# stamps that same style onto the cell too, whereas the authored row just
# inherits it from the column. Reset to Normal so the cell stays unstyled,
# matching the other data rows in this table.
$wsTypo.Range("B7:E7").Style = "Normal"
# G7/H7/J7 stay blank in the authored row, but are still present as empty
# cells (part of the table row). Touch them with a format no-op so they get
# materialised without picking up an explicit style.
$wsTypo.Range("G7").Font.Bold = $false
$wsTypo.Range("H7").Font.Bold = $false
$wsTypo.Range("J7").Font.Bold = $false

# ---- Translation sheet: row 62 default text changes from "Activity" to "Cycling" ----
$wsTrans.Range("F62").Value = "Cycling"

# ---- Translation sheet: new rows 63-84 for the activity user screen ----
$rows = @(
    @(63, "SingleUseId100", "Narrow", "Right",  "LTR", "<value>"),
    @(64, "SingleUseId101", "Narrow", "Left",   "LTR", "0.00"),
    @(65, "SingleUseId102", "Narrow", "Right",  "LTR", "<value>"),
    @(66, "SingleUseId103", "Narrow", "Left",   "LTR", "0.00"),
    @(67, "SingleUseId104", "Narrow", "Center", "LTR", "<value>"),
    @(68, "SingleUseId105", "Narrow", "Left",   "LTR", "00:00:00.0"),
    @(69, "SingleUseId106", "Small",  "Center", "LTR", "Timer"),
    @(70, "SingleUseId107", "Small",  "Center", "LTR", "Distance"),
    @(71, "SingleUseId108", "Small",  "Center", "LTR", "Speed"),
    @(72, "SingleUseId109", "Small",  "Center", "LTR", "Avg speed"),
    @(73, "SingleUseId110", "Small",  "Center", "LTR", "Max speed"),
    @(74, "SingleUseId111", "Narrow", "Right",  "LTR", "<value>"),
    @(75, "SingleUseId112", "Narrow", "Left",   "LTR", "0.00"),
    @(76, "SingleUseId113", "Small",  "Center", "LTR", "Clock"),
    @(77, "SingleUseId114", "Narrow", "Center", "LTR", "<value>"),
    @(78, "SingleUseId115", "Narrow", "Left",   "LTR", "12:00:00"),
    @(79, "SingleUseId98",  "Narrow", "Right",  "LTR", "<value>"),
    @(80, "SingleUseId117", "Small",  "Left",   "LTR", "km`n h"),
    @(81, "SingleUseId118", "Small",  "Left",   "LTR", "km`n h"),
    @(82, "SingleUseId119", "Narrow", "Left",   "LTR", "0.00"),
    @(83, "SingleUseId120", "Small",  "Left",   "LTR", "km`n h"),
    @(84, "SingleUseId121", "Small",  "Left",   "LTR", "km")
)

foreach ($row in $rows) {
    $r = $row[0]
    $wsTrans.Range("B$r").Value = $row[1]
    $wsTrans.Range("C$r").Value = $row[2]
    $wsTrans.Range("D$r").Value = $row[3]
    $wsTrans.Range("E$r").Value = $row[4]

    $default = $row[5]
    if ($default -eq "0.00") {
        Set-TextValue $wsTrans.Range("F$r") $default
    } else {
        $wsTrans.Range("F$r").Value = $default
    }
}
